# Update cryptos worksheet values to reflect the latest scrape
# (auto-generated from the OOXML diff of the commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.924.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.356.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.96%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.75"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.44"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.416"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.937.55"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.17"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.374.34"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.055.26"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.19"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.84"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.14%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.32"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.561"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.518.87"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.176"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.88"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.39"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "169.67"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.76"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "29.29"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.396.64"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0754"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.37"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.760"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.30"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.33%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.512.08"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.83"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.49%  "
